$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.681.53"
$ws.Range("E2").Value = "  -2.00%  "
$ws.Range("D3").Value = "1.758.58"
$ws.Range("E3").Value = "  -1.98%  "
$ws.Range("D4").Value = "'1.006"
$ws.Range("E4").Value = "  +0.25%  "
$ws.Range("D5").Value = "'326.31"
$ws.Range("E5").Value = "  -0.20%  "
$ws.Range("E6").Value = "  +0.27%  "
$ws.Range("D7").Value = "'0.4446"
$ws.Range("E7").Value = "  -0.84%  "
$ws.Range("D8").Value = "'0.3727"
$ws.Range("E8").Value = "  -0.62%  "
$ws.Range("D9").Value = "'45.60"
$ws.Range("E9").Value = "  +1.97%  "
$ws.Range("D10").Value = "'0.07791"
$ws.Range("E10").Value = "  +3.58%  "
$ws.Range("D11").Value = "'1.127"
$ws.Range("E11").Value = "  -1.55%  "
$ws.Range("D12").Value = "'1.003"
$ws.Range("E12").Value = "  +0.11%  "
$ws.Range("D13").Value = "'21.78"
$ws.Range("E13").Value = "  -3.51%  "
$ws.Range("D14").Value = "'6.202"
$ws.Range("E14").Value = "  -1.27%  "
$ws.Range("D15").Value = "'7.381"
$ws.Range("E15").Value = "  -2.06%  "
$ws.Range("D16").Value = "1.759.71"
$ws.Range("E16").Value = "  -1.84%  "
$ws.Range("D17").Value = "'91.27"
$ws.Range("E17").Value = "  +12.66%  "
$ws.Range("E18").Value = "  -0.64%  "
$ws.Range("D19").Value = "'0.06251"
$ws.Range("E19").Value = "  -7.09%  "
$ws.Range("E20").Value = "  +0.24%  "
$ws.Range("E21").Value = "  -0.72%  "
$ws.Range("D22").Value = "'6.193"
$ws.Range("E22").Value = "  -2.16%  "
$ws.Range("D23").Value = "'0.5327"
$ws.Range("E23").Value = "  -3.15%  "
$ws.Range("D24").Value = "27.713.52"
$ws.Range("E24").Value = "  -1.84%  "
$ws.Range("D25").Value = "'11.67"
$ws.Range("E25").Value = "  -0.65%  "
$ws.Range("D26").Value = "'2.332"
$ws.Range("E26").Value = "  -3.67%  "
$ws.Range("D27").Value = "'20.84"
$ws.Range("E27").Value = "  +1.91%  "
$ws.Range("D28").Value = "'153.59"
$ws.Range("E28").Value = "  +0.99%  "
$ws.Range("D29").Value = "'2.352"
$ws.Range("E29").Value = "  -0.09%  "
$ws.Range("D30").Value = "1.958.88"
$ws.Range("E30").Value = "  -1.46%  "
$ws.Range("D31").Value = "'129.10"
$ws.Range("E31").Value = "  -2.70%  "
$ws.Range("D32").Value = "'1.214"
$ws.Range("E32").Value = "  -0.99%  "
$ws.Range("D33").Value = "'5.777"
$ws.Range("E33").Value = "  -0.49%  "
$ws.Range("D34").Value = "'0.09258"
$ws.Range("E34").Value = "  -1.64%  "
$ws.Range("D35").Value = "'3.692"
$ws.Range("E35").Value = "  -8.23%  "
$ws.Range("D36").Value = "'12.77"
$ws.Range("E36").Value = "  +5.54%  "
$ws.Range("D37").Value = "'0.02346"
$ws.Range("E37").Value = "  +0.84%  "
$ws.Range("D38").Value = "'0.2190"
$ws.Range("E38").Value = "  -5.86%  "
$ws.Range("D39").Value = "'0.6504"
$ws.Range("E39").Value = "  -0.67%  "
$ws.Range("D40").Value = "'5.092"
$ws.Range("E40").Value = "  -1.42%  "
$ws.Range("E41").Value = "  -3.26%  "
$ws.Range("D42").Value = "'1.193"
$ws.Range("E42").Value = "  -1.17%  "
$ws.Range("D43").Value = "'8.026"
$ws.Range("E43").Value = "  -3.39%  "
$ws.Range("B44").Value = "Frax"
$ws.Range("C44").Value = "https://coinranking.com/coin/KfWtaeV1W+frax-frax"
$ws.Range("D44").Value = "'1.003"
$ws.Range("E44").Value = "  +0.26%  "
$ws.Range("B45").Value = "WEMIXTOKEN"
$ws.Range("C45").Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Range("D45").Value = "'1.414"
$ws.Range("E45").Value = "  -3.74%  "
$ws.Range("D46").Value = "'13.88"
$ws.Range("E46").Value = "  -1.41%  "
$ws.Range("D47").Value = "'0.6009"
$ws.Range("E47").Value = "  -1.33%  "
$ws.Range("D48").Value = "'3.751"
$ws.Range("D49").Value = "'125.88"
$ws.Range("E49").Value = "  -3.10%  "
$ws.Range("E50").Value = "  -1.15%  "
$ws.Range("D51").Value = "'1.147"
$ws.Range("E51").Value = "  -1.07%  "
